$wb = $excel.ActiveWorkbook

# --- Rename sheets (task order ids) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16511687794108207"
$wb.Worksheets.Item(2).Name = "NB_TO-1651168780956945"
$wb.Worksheets.Item(3).Name = "RS_TO-16511687809579477"
$wb.Worksheets.Item(4).Name = "TOL_TO-1651168781015208"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511687810926366"

# --- Sheet 1 (GNG) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1651168779367661.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687793949807.csv"
$ws1.Range("B4").Value = "go_stims-16511687793959787.csv"
$ws1.Range("B5").Value = "GNG_stims-1651168779409823.csv"

# --- Sheet 2 (NB) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16511687806125712.csv"
$ws2.Range("B3").Value = "ZB-match_0-1651168779997432.csv"
$ws2.Range("B4").Value = "OB-1651168780377066.csv"
$ws2.Range("B5").Value = "ZB-match_0-16511687797617807.csv"
$ws2.Range("B6").Value = "TB-1651168780906191.csv"
$ws2.Range("B7").Value = "OB-16511687804901493.csv"
$ws2.Range("B8").Value = "TB-16511687809443426.csv"
$ws2.Range("B9").Value = "OB-16511687801603131.csv"
$ws2.Range("B10").Value = "ZB-match_5-16511687795500922.csv"

# --- Sheet 3 (RS) --- (no cell content changes, only sheet name changed above)

# --- Sheet 4 (TOL) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511687809721906.csv"
$ws4.Range("B3").Value = "ZM_stims-16511687809589477.csv"
$ws4.Range("B4").Value = "MM_stims-16511687809992924.csv"
$ws4.Range("B5").Value = "ZM_stims-1651168780973162.csv"
$ws4.Range("B6").Value = "MM_stims-16511687810142097.csv"
$ws4.Range("B7").Value = "ZM_stims-16511687810002856.csv"

# --- Sheet 5 (vSAT) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16511687810219784.csv"
$ws5.Range("B3").Value = "vSAT_stims-1651168781061874.csv"
$ws5.Range("B4").Value = "vSAT_stims-16511687810782502.csv"
$ws5.Range("B5").Value = "SAT_stims-1651168781046031.csv"
